# Refresh the cryptos price/volume table (Sheet1) to the latest scrape.
# Coin name / link cells only change where the underlying GitHub Actions
# scrape reordered two neighbouring rows (rows 39/40, 41/42, 43/44 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.851.16'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.63%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.253.81'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.97%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.89'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.80%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.17'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.28%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.23%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.134'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.32%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.43%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.415'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.01%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.812.97'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.78%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.137'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.64%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.67'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.00%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.857.89'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.64%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.04%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.250.56'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.91%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.83'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.39%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.53'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.65%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '378.99'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.64%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.64'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.36%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.10%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.27'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.22%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.513'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000120'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.27%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.97'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.20%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.58%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.08%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.98'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.54%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.68'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.78%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.86'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.01'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.62%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.27'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.62%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.49'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.33%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.845'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.20%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.87'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.80%  '

# Row 39
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.46'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.63%  '

# Row 40
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.60'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.47%  '

# Row 41
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.34%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.59'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.88%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.68'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.05%  '

# Row 44
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '346.95'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.17%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.08'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.92%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0686'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.48%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.635.25'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.52%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0286'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.35%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.992'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.19%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.16'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.08%  '
